$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.34808886051178
$ws.Range("B1").Value = 2.48285174369812
$ws.Range("C1").Value = 4.926457405090332
$ws.Range("D1").Value = 2.364908695220947
$ws.Range("E1").Value = 0.913968563079834
